# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worksheet lists mora periods per worker. Previously the periods for
# the two main workers (LUIS RODOLFO CAMPUZANO DE AVILA / CC 1047464119 and
# JULIO ENRIQUE CASTELLON AVENDAÑO / CC 1002476864) were interleaved row by
# row; the old "2301" period didn't exist yet. This edit replaces that
# block (rows 16-42) with the updated data: each worker's 13 periods
# (2301 down to 2201) grouped together, plus a corrected "Valor Mora"
# (column G) for JONH FREDY PIRATEQUE PLAZAS's existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Descending list of periods, newest ("2301") first, then 2212 .. 2201.
$periods = @("2301","2212","2211","2210","2209","2208","2207","2206","2205","2204","2203","2202","2201")

# Column F (Salario Basico date-serial column): the "2301" row uses 30284,
# every other period row uses 36341 — matches the source data exactly.
function Get-FValue($period) {
    if ($period -eq "2301") { return 30284 }
    return 36341
}

# Rows 16-28: LUIS RODOLFO CAMPUZANO DE AVILA, one row per period.
$startRow = 16
for ($i = 0; $i -lt $periods.Count; $i++) {
    $r = $startRow + $i
    $period = $periods[$i]
    $ws.Range("C$r").Value = "1047464119"
    $ws.Range("D$r").Value = "LUIS RODOLFO CAMPUZANO DE AVILA"
    $ws.Range("E$r").Value = $period
    $ws.Range("F$r").Value = Get-FValue $period
    $ws.Range("G$r").Value = 908526
}

# Row 29: JONH FREDY PIRATEQUE PLAZAS, period 2207 — Valor Mora corrected.
$ws.Range("C29").Value = "88208149"
$ws.Range("D29").Value = "JONH FREDY PIRATEQUE PLAZAS"
$ws.Range("E29").Value = "2207"
$ws.Range("F29").Value = 40000
$ws.Range("G29").Value = 781242

# Rows 30-42: JULIO ENRIQUE CASTELLON AVENDAÑO, one row per period.
$startRow = 30
for ($i = 0; $i -lt $periods.Count; $i++) {
    $r = $startRow + $i
    $period = $periods[$i]
    $ws.Range("C$r").Value = "1002476864"
    $ws.Range("D$r").Value = "JULIO ENRIQUE CASTELLON AVENDAÑO"
    $ws.Range("E$r").Value = $period
    $ws.Range("F$r").Value = Get-FValue $period
    $ws.Range("G$r").Value = 908526
}
